$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G3").Value = 1.8
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 5.25
$ws.Range("J3").Value = 2.6
$ws.Range("L3").Value = 6
$ws.Range("O3").Value = 1.57
$ws.Range("P3").Value = 2.38
$ws.Range("Q3").Value = 2.7
$ws.Range("R3").Value = 1.44
$ws.Range("U3").Value = 2.5
$ws.Range("V3").Value = 1.5
$ws.Range("W3").Value = 4.75
$ws.Range("X3").Value = 6.5
$ws.Range("Y3").Value = 9.5
$ws.Range("Z3").Value = 13
$ws.Range("AA3").Value = 19
$ws.Range("AD3").Value = 7
$ws.Range("AH3").Value = 9
$ws.Range("AI3").Value = 23
$ws.Range("AJ3").Value = 19
$ws.Range("AL3").Value = 51
$ws.Range("AM3").Value = 67
$ws.Range("AN3").Value = 3.5
$ws.Range("AO3").Value = 10
$ws.Range("AW3").Value = 6.5
$ws.Range("AX3").Value = 34
$ws.Range("AZ3").Value = 151
$ws.Range("BA3").Value = 201
$ws.Range("G4").Value = 1.75
$ws.Range("I4").Value = 6.25
$ws.Range("K4").Value = 1.91
$ws.Range("Q4").Value = 2.88
$ws.Range("R4").Value = 1.4
$ws.Range("AN4").Value = 3.4
$ws.Range("M6").Value = 1.01
$ws.Range("N6").Value = 14.8
$ws.Range("O6").Value = 1.14
$ws.Range("P6").Value = 5.24
$ws.Range("S6").Value = 1.24
$ws.Range("T6").Value = 3.93
$ws.Range("U6").Value = 1.69
$ws.Range("V6").Value = 2.1
$ws.Range("BC8").Value = 126
$ws.Range("M9").Value = 1.02
$ws.Range("O9").Value = 1.11
$ws.Range("M10").Value = 1.05
$ws.Range("O10").Value = 1.37
$ws.Range("Q10").Value = 2.35
$ws.Range("R10").Value = 1.57
$ws.Range("G11").Value = 2
$ws.Range("I11").Value = 4.1
$ws.Range("J11").Value = 2.88
$ws.Range("O11").Value = 1.54
$ws.Range("P11").Value = 2.25
$ws.Range("W11").Value = 5.5
$ws.Range("X11").Value = 8
$ws.Range("AE11").Value = 21
$ws.Range("AF11").Value = 81
$ws.Range("AK11").Value = 41
$ws.Range("M12").Value = 1.08
$ws.Range("O12").Value = 1.44
$ws.Range("P12").Value = 2.63
$ws.Range("G13").Value = 1.27
$ws.Range("H13").Value = 5
$ws.Range("I13").Value = 12
$ws.Range("L13").Value = 9
$ws.Range("N13").Value = 15
$ws.Range("O13").Value = 1.2
$ws.Range("Q13").Value = 1.65
$ws.Range("R13").Value = 2.2
$ws.Range("X13").Value = 6
$ws.Range("Z13").Value = 7.5
$ws.Range("AD13").Value = 10
$ws.Range("AE13").Value = 23
$ws.Range("AH13").Value = 26
$ws.Range("AJ13").Value = 34
$ws.Range("AK13").Value = 151
$ws.Range("AL13").Value = 81
